# Festsetzung_Sammelbescheid.docx edits:
#  1) Move the _GoBack bookmark from its old location (end of the
#     "Art = 13 ..." paragraph) to the very start of the document
#     (first, empty paragraph) - a Word "last edit position" artifact.
#  2) Resize several columns of the big table (the one tracking the
#     individual Teilgewaesserbenutzungen) so the "erlaubte Menge"
#     column widens and neighbouring columns shrink to compensate.

$d = $word.ActiveDocument

# --- 1) Move the _GoBack bookmark -----------------------------------
$goBack = $null
try {
    $goBack = $d.Bookmarks.Item("_GoBack")
} catch {
    $goBack = $null
}
if ($goBack -ne $null) {
    $goBack.Delete()
}
$firstParaStart = $d.Paragraphs.Item(1).Range.Start
$bookmarkRange = $d.Range($firstParaStart, $firstParaStart)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- 2) Resize table columns -----------------------------------------
# The table is the big 12-column / 4-row table further down the
# document (the grid listing Wasserbenutzung / Art / Menge / ...).
$table = $d.Tables.Item(2)

# width values are in dxa (twentieths of a point); Word's Width
# property is expressed in points, so divide by 20.
$table.Cell(2, 1).Width  = 1242 / 20.0   # 568  -> 1242
$table.Cell(2, 2).Width  = 1134 / 20.0   # 1667 -> 1134
$table.Cell(2, 3).Width  = 709  / 20.0   # 850  -> 709
$table.Cell(2, 6).Width  = 1560 / 20.0   # 1276 -> 1560
$table.Cell(2, 7).Width  = 708  / 20.0   # 992  -> 708
